# Update cryptocurrency price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.060.05'
$ws.Range('E2').Value = '  +3.65%  '
$ws.Range('D3').Value = '1.726.26'
$ws.Range('E3').Value = '  +2.96%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = "'218.84"
$ws.Range('E5').Value = '  +1.61%  '
$ws.Range('E6').Value = '  +1.11%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('E8').Value = '  +13.14%  '
$ws.Range('E9').Value = '  +3.15%  '
$ws.Range('E10').Value = '  +1.81%  '
$ws.Range('E11').Value = '  +1.77%  '
$ws.Range('D12').Value = '1.971.06'
$ws.Range('E12').Value = '  +3.05%  '
$ws.Range('D13').Value = '1.735.48'
$ws.Range('E13').Value = '  +3.27%  '
$ws.Range('D14').Value = "'4.25"
$ws.Range('E14').Value = '  +3.27%  '
$ws.Range('D15').Value = "'0.562"
$ws.Range('E15').Value = '  +5.34%  '
$ws.Range('D16').Value = "'67.56"
$ws.Range('E16').Value = '  +2.38%  '
$ws.Range('D17').Value = '28.019.58'
$ws.Range('E17').Value = '  +3.54%  '
$ws.Range('D18').Value = "'242.65"
$ws.Range('E18').Value = '  +2.17%  '
$ws.Range('E19').Value = '  +1.80%  '
$ws.Range('E20').Value = '  -3.23%  '
$ws.Range('E21').Value = '  -0.15%  '
$ws.Range('D22').Value = "'4.62"
$ws.Range('E22').Value = '  +3.66%  '
$ws.Range('E23').Value = '  +4.25%  '
$ws.Range('E24').Value = '  -0.08%  '
$ws.Range('D25').Value = "'148.93"
$ws.Range('E25').Value = '  +1.45%  '
$ws.Range('E26').Value = '  +4.31%  '
$ws.Range('E27').Value = '  +2.46%  '
$ws.Range('E28').Value = '  +1.58%  '
$ws.Range('E29').Value = '  +0.10%  '
$ws.Range('E30').Value = '  +2.45%  '
$ws.Range('E31').Value = '  +2.28%  '
$ws.Range('E32').Value = '  +2.72%  '
$ws.Range('D33').Value = '1.494.03'
$ws.Range('E33').Value = '  -3.35%  '
$ws.Range('D34').Value = "'3.26"
$ws.Range('E34').Value = '  +2.43%  '
$ws.Range('D35').Value = "'1.65"
$ws.Range('E35').Value = '  -2.54%  '
$ws.Range('D36').Value = "'0.952"
$ws.Range('E36').Value = '  +3.00%  '
$ws.Range('E37').Value = '  +1.09%  '
$ws.Range('E38').Value = '  +0.66%  '
$ws.Range('E39').Value = '  +0.44%  '
$ws.Range('E40').Value = '  +0.16%  '
$ws.Range('E41').Value = '  +4.73%  '
$ws.Range('D42').Value = "'5.84"
$ws.Range('E42').Value = '  +4.19%  '
$ws.Range('E43').Value = '  -0.10%  '
$ws.Range('E44').Value = '  +2.37%  '
$ws.Range('D45').Value = '1.874.72'
$ws.Range('E45').Value = '  +2.79%  '
$ws.Range('E46').Value = '  +1.96%  '
$ws.Range('E47').Value = '  +11.79%  '
$ws.Range('D48').Value = "'91.11"
$ws.Range('E48').Value = '  +0.49%  '
$ws.Range('E49').Value = '  +4.08%  '
$ws.Range('E50').Value = '  +0.79%  '
$ws.Range('D51').Value = "'8.20"
$ws.Range('E51').Value = '  +2.14%  '
